$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(4, "Estados Unidos", 644348, 259, 48708, 567086, 13487, 25, 28554),
    @(5, "España", 180659, 0, 70853, 90994, 7371, 0, 18812),
    @(6, "Italia", 165155, 0, 38092, 105418, 3079, 0, 21645),
    @(7, "Francia", 147863, 0, 30955, 99741, 6457, 0, 17167),
    @(8, "Alemania", 134753, 0, 77000, 53949, 4288, 0, 3804),
    @(9, "Reino Unido", 98476, 0, 0, 85264, 1559, 0, 12868),
    @(10, "China", 82341, 46, 77892, 1107, 95, 0, 3342),
    @(11, "Iran", 76389, 0, 49933, 21679, 3643, 0, 4777),
    @(12, "Turquia", 69392, 0, 5674, 62200, 1820, 0, 1518),
    @(13, "Belgica", 34809, 1236, 7562, 22390, 1182, 417, 4857),
    @(14, "Brasil", 28912, 302, 14026, 13126, 296, 3, 1760),
    @(15, "Canada", 28379, 0, 8979, 18390, 557, 0, 1010),
    @(16, "Paises Bajos", 28153, 0, 250, 24769, 1279, 0, 3134),
    @(17, "Rusia", 27938, 3448, 2304, 25402, 8, 34, 232),
    @(18, "Suiza", 26336, 0, 15400, 9687, 386, 10, 1249),
    @(19, "Portugal", 18091, 0, 383, 17109, 208, 0, 599),
    @(20, "Austria", 14412, 62, 8098, 5921, 232, 0, 393),
    @(21, "Israel", 12591, 90, 2624, 9827, 174, 10, 140),
    @(22, "Irlanda", 12547, 0, 77, 12026, 158, 0, 444),
    @(23, "India", 12456, 86, 1513, 10520, 0, 1, 423),
    @(24, "Suecia", 11927, 0, 381, 10343, 954, 0, 1203),
    @(25, "Peru", 11475, 0, 3108, 8113, 146, 0, 254),
    @(26, "Corea del Sur", 10613, 22, 7757, 2627, 55, 4, 229),
    @(27, "Japon", 8626, 0, 901, 7547, 168, 0, 178),
    @(28, "Chile", 8273, 0, 2937, 5242, 389, 0, 94),
    @(29, "Ecuador", 7858, 0, 780, 6690, 135, 0, 388),
    @(30, "Polonia", 7771, 189, 774, 6705, 160, 6, 292),
    @(31, "Rumania", 7216, 0, 1217, 5612, 245, 15, 387),
    @(32, "Dinamarca", 6879, 198, 2748, 3822, 89, 0, 309),
    @(33, "Noruega", 6798, 1, 32, 6616, 64, 0, 150),
    @(34, "Pakistan", 6505, 122, 1645, 4736, 46, 13, 124),
    @(35, "Australia", 6468, 21, 3747, 2658, 66, 0, 63),
    @(36, "Chequia", 6303, 2, 831, 5306, 75, 0, 166),
    @(37, "Arabia Saudita", 5862, 0, 931, 4852, 71, 0, 79),
    @(38, "Mexico", 5847, 448, 2125, 3273, 207, 43, 449),
    @(39, "Filipinas", 5660, 207, 435, 4863, 1, 13, 362),
    @(40, "Indonesia", 5516, 380, 548, 4470, 0, 29, 498),
    @(41, "Emiratos Arabes Unidos", 5365, 0, 1034, 4298, 1, 0, 33),
    @(42, "Malasia", 5182, 110, 2766, 2332, 56, 1, 84),
    @(43, "Serbia", 4873, 0, 400, 4374, 128, 0, 99),
    @(44, "Ucrania", 4161, 397, 186, 3859, 45, 8, 116),
    @(45, "Panama", 3751, 0, 75, 3573, 106, 0, 103),
    @(46, "Bielorrusia", 3728, 0, 203, 3489, 68, 0, 36),
    @(47, "Catar", 3711, 0, 406, 3298, 37, 0, 7),
    @(48, "Singapur", 3699, 0, 652, 3037, 29, 0, 10),
    @(49, "Republica Dominicana", 3614, 0, 208, 3217, 143, 0, 189),
    @(50, "Luxemburgo", 3373, 0, 526, 2778, 33, 0, 69),
    @(51, "Finlandia", 3369, 132, 300, 2997, 75, 0, 72),
    @(52, "Colombia", 3105, 0, 452, 2522, 106, 0, 131),
    @(53, "Tailandia", 2672, 29, 1593, 1033, 61, 3, 46),
    @(54, "Argentina", 2571, 0, 596, 1863, 117, 0, 112),
    @(55, "Sudafrica", 2506, 0, 410, 2062, 7, 0, 34),
    @(56, "Egipto", 2505, 0, 589, 1733, 0, 0, 183),
    @(57, "Marruecos", 2251, 227, 247, 1876, 1, 1, 128),
    @(58, "Grecia", 2192, 0, 269, 1821, 72, 0, 102),
    @(59, "Argelia", 2160, 0, 708, 1116, 60, 0, 336),
    @(60, "Moldavia", 2049, 0, 235, 1763, 80, 5, 51),
    @(61, "Croacia", 1741, 0, 473, 1235, 31, 0, 33),
    @(62, "Islandia", 1727, 0, 1077, 642, 8, 0, 8),
    @(63, "Barein", 1673, 2, 663, 1003, 3, 0, 7),
    @(64, "Hungria", 1652, 73, 199, 1311, 58, 8, 142),
    @(65, "Banglades", 1572, 341, 49, 1463, 1, 10, 60),
    @(66, "Kuwait", 1524, 119, 225, 1296, 32, 0, 3),
    @(67, "Estonia", 1434, 34, 133, 1265, 10, 1, 36),
    @(68, "Irak", 1415, 0, 812, 524, 0, 0, 79),
    @(69, "Nueva Zelanda", 1401, 15, 770, 622, 3, 0, 9),
    @(70, "Uzbekistan", 1349, 47, 107, 1238, 8, 0, 4),
    @(71, "Kazajistan", 1341, 46, 263, 1061, 22, 1, 17),
    @(72, "Azerbaiyan", 1253, 0, 404, 836, 24, 0, 13),
    @(73, "Eslovenia", 1248, 0, 165, 1022, 34, 0, 61),
    @(74, "Armenia", 1159, 48, 358, 783, 30, 1, 18),
    @(75, "Lituania", 1128, 37, 178, 920, 14, 0, 30),
    @(76, "Bosnia y Herzegovina", 1110, 0, 253, 816, 4, 0, 41),
    @(77, "Oman", 1019, 109, 176, 839, 3, 0, 4),
    @(78, "Hong Kong", 1018, 1, 485, 529, 9, 0, 4),
    @(79, "Republica de Macedonia", 974, 0, 98, 831, 15, 0, 45),
    @(80, "Eslovaquia", 863, 0, 151, 706, 5, 0, 6),
    @(81, "Camerun", 848, 0, 165, 666, 0, 0, 17),
    @(82, "Cuba", 814, 0, 151, 639, 15, 0, 24),
    @(83, "Afganistan", 784, 0, 43, 716, 0, 0, 25),
    @(84, "Bulgaria", 783, 36, 122, 624, 31, 1, 37),
    @(85, "Tunez", 780, 0, 43, 702, 89, 0, 35),
    @(86, "Republica de Chipre", 715, 0, 65, 638, 8, 0, 12),
    @(87, "Crucero", 712, 0, 644, 56, 7, 0, 12),
    @(88, "Letonia", 675, 9, 57, 613, 3, 0, 5),
    @(89, "Principado de Andorra", 673, 0, 169, 471, 17, 0, 33),
    @(90, "Libano", 658, 0, 85, 552, 30, 0, 21),
    @(91, "Costa de Marfil", 654, 0, 146, 502, 0, 0, 6),
    @(92, "Ghana", 641, 0, 83, 550, 2, 0, 8),
    @(93, "Costa Rica", 626, 0, 67, 555, 11, 0, 4),
    @(94, "Niger", 584, 0, 90, 480, 0, 0, 14),
    @(95, "Burkina Faso", 542, 0, 226, 284, 0, 0, 32),
    @(96, "Albania", 494, 0, 251, 218, 5, 0, 25),
    @(97, "Uruguay", 493, 0, 272, 212, 13, 0, 9),
    @(98, "Kirguistan", 466, 17, 91, 370, 5, 0, 5),
    @(99, "Bolivia", 441, 44, 14, 398, 3, 1, 29),
    @(100, "Republica de Yibuti", 435, 0, 71, 362, 0, 0, 2),
    @(101, "Honduras", 426, 7, 9, 382, 10, 4, 35),
    @(102, "Nigeria", 407, 0, 128, 267, 2, 0, 12),
    @(103, "Guinea", 404, 0, 31, 372, 0, 0, 1),
    @(104, "Jordania", 401, 0, 250, 144, 5, 0, 7),
    @(105, "Malta", 399, 0, 82, 314, 4, 0, 3),
    @(106, "Taiwan", 395, 0, 155, 234, 0, 0, 6),
    @(107, "San Marino", 393, 0, 53, 304, 15, 0, 36),
    @(108, "Reunion", 391, 0, 237, 154, 3, 0, 0),
    @(109, "Estado de Palestina", 374, 0, 63, 309, 0, 0, 2),
    @(110, "Georgia", 336, 30, 74, 259, 6, 0, 3),
    @(111, "Mauricio", 324, 0, 65, 250, 3, 0, 9),
    @(112, "Senegal", 314, 0, 190, 122, 1, 0, 2),
    @(113, "Montenegro", 290, 2, 55, 231, 7, 0, 4),
    @(114, "Isla de Man", 283, 27, 153, 126, 13, 0, 4),
    @(115, "Vietnam", 268, 0, 171, 97, 8, 0, 0),
    @(116, "Consejo Danes para los Refugiados", 267, 13, 23, 222, 0, 1, 22),
    @(117, "Sri Lanka", 238, 0, 65, 166, 1, 0, 7),
    @(118, "Kenia", 225, 0, 53, 162, 2, 0, 10),
    @(119, "Mayotte", 217, 0, 69, 145, 3, 0, 3),
    @(120, "Venezuela", 197, 0, 111, 77, 6, 0, 9),
    @(121, "Guatemala", 196, 16, 19, 172, 3, 0, 5),
    @(122, "Islas Feroe", 184, 0, 169, 15, 0, 0, 0),
    @(123, "Paraguay", 174, 13, 30, 136, 1, 0, 8),
    @(124, "El Salvador", 164, 5, 33, 125, 2, 0, 6),
    @(125, "Martinica", 158, 0, 73, 77, 17, 0, 8),
    @(126, "Mali", 148, 0, 34, 101, 0, 0, 13),
    @(127, "Guadalupe", 145, 0, 67, 70, 13, 0, 8),
    @(128, "Ruanda", 136, 0, 54, 82, 0, 0, 0),
    @(129, "Brunei", 136, 0, 108, 27, 2, 0, 1),
    @(130, "Gibraltar", 131, 0, 104, 27, 1, 0, 0),
    @(131, "Jamaica", 125, 0, 21, 99, 0, 0, 5),
    @(132, "Camboya", 122, 0, 98, 24, 1, 0, 0),
    @(133, "Congo", 117, 0, 11, 101, 0, 0, 5),
    @(134, "Trinidad yTobago", 114, 0, 20, 86, 0, 0, 8),
    @(135, "Madagascar", 110, 0, 29, 81, 1, 0, 0),
    @(136, "Monaco", 93, 0, 12, 78, 2, 0, 3),
    @(137, "Aruba", 93, 0, 39, 53, 1, 0, 1),
    @(138, "Etiopia", 92, 7, 15, 74, 0, 0, 3),
    @(139, "Tanzania", 88, 0, 11, 73, 0, 0, 4),
    @(140, "Guayana Francesa", 86, 0, 51, 35, 1, 0, 0),
    @(141, "Birmania", 85, 11, 2, 79, 0, 0, 4),
    @(142, "Bermudas", 81, 0, 33, 43, 3, 0, 5),
    @(143, "Togo", 81, 0, 35, 43, 0, 0, 3),
    @(144, "Gabon", 80, 0, 4, 75, 0, 0, 1),
    @(145, "Somalia", 80, 0, 2, 73, 2, 0, 5),
    @(146, "Liechtenstein", 79, 0, 55, 23, 0, 0, 1),
    @(147, "Barbados", 75, 2, 15, 55, 4, 0, 5),
    @(148, "Islas Caimanes", 60, 0, 6, 53, 3, 0, 1),
    @(149, "Liberia", 59, 0, 4, 49, 0, 0, 6),
    @(150, "Cabo Verde", 56, 0, 1, 54, 0, 0, 1),
    @(151, "Polinesia Francesa", 55, 0, 0, 55, 1, 0, 0),
    @(152, "Uganda", 55, 0, 12, 43, 0, 0, 0),
    @(153, "Guyana", 55, 0, 8, 41, 5, 0, 6),
    @(154, "San Martin (Parte Holandesa)", 53, 0, 5, 39, 2, 0, 9),
    @(155, "Bahamas", 53, 0, 6, 39, 1, 0, 8),
    @(156, "Guinea Ecuatorial", 51, 0, 4, 47, 0, 0, 0),
    @(157, "Libia", 48, 0, 11, 36, 0, 0, 1),
    @(158, "Zambia", 48, 0, 30, 16, 1, 0, 2),
    @(159, "Macao", 45, 0, 16, 29, 1, 0, 0),
    @(160, "Guinea-Bisau", 43, 0, 0, 43, 0, 0, 0),
    @(161, "Haiti", 41, 0, 0, 38, 0, 0, 3),
    @(162, "Puerto Rico", 39, 0, 1, 36, 0, 0, 2),
    @(163, "Eritrea", 35, 0, 0, 35, 0, 0, 0),
    @(164, "San Martin (Parte Francesa)", 35, 0, 13, 20, 5, 0, 2),
    @(165, "Benin", 35, 0, 18, 16, 0, 0, 1),
    @(166, "Siria", 33, 0, 5, 26, 0, 0, 2),
    @(167, "Guam", 32, 0, 0, 31, 0, 0, 1),
    @(168, "Sudan", 32, 0, 4, 23, 0, 0, 5),
    @(169, "Mongolia", 31, 1, 5, 26, 0, 0, 0),
    @(170, "Mozambique", 29, 0, 2, 27, 0, 0, 0),
    @(171, "Republica del Chad", 27, 4, 5, 22, 0, 0, 0),
    @(172, "Zimbabue", 23, 0, 1, 19, 0, 0, 3),
    @(173, "Antigua y Barbuda", 23, 0, 3, 17, 1, 1, 3),
    @(174, "Maldivas", 23, 1, 16, 7, 0, 0, 0),
    @(175, "Laos", 19, 0, 2, 17, 0, 0, 0),
    @(176, "Angola", 19, 0, 5, 12, 0, 0, 2),
    @(177, "Nueva Caledonia", 18, 0, 1, 17, 1, 0, 0),
    @(178, "Timor Oriental", 18, 10, 1, 17, 0, 0, 0),
    @(179, "Belice", 18, 0, 0, 16, 1, 0, 2),
    @(180, "Islas Virgenes de los Estados Unidos", 17, 0, 0, 17, 0, 0, 0),
    @(181, "Fiyi", 17, 1, 0, 17, 0, 0, 0),
    @(182, "Nepal", 16, 0, 1, 15, 0, 0, 0),
    @(183, "Malaui", 16, 0, 0, 14, 1, 0, 2),
    @(184, "Namibia", 16, 0, 3, 13, 0, 0, 0),
    @(185, "Dominica", 16, 0, 8, 8, 0, 0, 0),
    @(186, "Suazilandia", 16, 1, 8, 8, 0, 0, 0),
    @(187, "Botsuana", 15, 2, 0, 14, 0, 0, 1),
    @(188, "Santa Lucia", 15, 0, 11, 4, 0, 0, 0),
    @(189, "Granada", 14, 0, 0, 14, 2, 0, 0),
    @(190, "San Cristobal y Nieves", 14, 0, 0, 14, 0, 0, 0),
    @(191, "Curazao", 14, 0, 10, 3, 0, 0, 1),
    @(192, "Sierra Leona", 13, 0, 0, 13, 0, 0, 0),
    @(193, "San Vicente y las Granadinas", 12, 0, 1, 11, 0, 0, 0),
    @(194, "Republica de Africa Central", 12, 0, 4, 8, 0, 0, 0),
    @(195, "Seychelles", 11, 0, 0, 11, 0, 0, 0),
    @(196, "Montserrat", 11, 0, 1, 10, 1, 0, 0),
    @(197, "Islas Malvinas", 11, 0, 1, 10, 0, 0, 0),
    @(198, "Groenlandia", 11, 0, 11, 0, 0, 0, 0),
    @(199, "Islas Turcas y Caicos", 10, 0, 0, 9, 0, 0, 1),
    @(200, "Surinam", 10, 0, 6, 3, 0, 0, 1),
    @(201, "Gambia", 9, 0, 2, 6, 0, 0, 1),
    @(202, "Nicaragua", 9, 0, 4, 4, 0, 0, 1),
    @(203, "Santa Sede", 8, 0, 2, 6, 0, 0, 0),
    @(204, "Mauritania", 7, 0, 2, 4, 0, 0, 1),
    @(205, "Sahara Occidental", 6, 0, 0, 6, 0, 0, 0),
    @(206, "San Bartolome", 6, 0, 4, 2, 0, 0, 0),
    @(207, "Burundi", 5, 0, 0, 4, 0, 0, 1),
    @(208, "Butan", 5, 0, 2, 3, 0, 0, 0),
    @(209, "Sudan del Sur", 4, 0, 0, 4, 0, 0, 0),
    @(210, "Santo Tome y Principe", 4, 0, 0, 4, 0, 0, 0),
    @(211, "Bonaire, San Eustaquio y Saba", 3, 0, 0, 3, 0, 0, 0),
    @(212, "Anguila", 3, 0, 1, 2, 0, 0, 0),
    @(213, "Islas Virgenes Britanicas", 3, 0, 2, 1, 0, 0, 0),
    @(214, "Papua Nueva Guinea", 2, 0, 0, 2, 0, 0, 0),
    @(215, "San Pedro y Miquelon", 1, 0, 0, 1, 0, 0, 0),
    @(216, "Yemen", 1, 0, 0, 1, 0, 0, 0)
)

foreach ($r in $data) {
    $rowNum = $r[0]
    $ws.Cells.Item($rowNum, 1).Value = $r[1]
    $ws.Cells.Item($rowNum, 2).Value = $r[2]
    $ws.Cells.Item($rowNum, 3).Value = $r[3]
    $ws.Cells.Item($rowNum, 4).Value = $r[4]
    $ws.Cells.Item($rowNum, 5).Value = $r[5]
    $ws.Cells.Item($rowNum, 6).Value = $r[6]
    $ws.Cells.Item($rowNum, 7).Value = $r[7]
    $ws.Cells.Item($rowNum, 8).Value = $r[8]
}

Write-Host "done"
